$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H28").Value = 758.51514
$ws.Range("I28").Value = 711.1905
$ws.Range("J28").Value = 841.3333
$ws.Range("K28").Value = 711.1905
$ws.Range("L28").Value = 841.3333
$ws.Range("M28").Value = -226.1905
$ws.Range("N28").Value = -1811.3333
$ws.Range("H98").Value = 2531.4546
$ws.Range("I98").Value = 2579.3125
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 2579.3125
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -1081.3125
$ws.Range("N98").Value = -3996
$ws.Range("H106").Value = 2132.6
$ws.Range("I106").Value = 2132.6
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2132.6
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -1501.6
$ws.Range("H122").Value = 2531.4546
$ws.Range("I122").Value = 2579.3125
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 7737.9375
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -5287.9375
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 2278.6
$ws.Range("I132").Value = 2494.9656
$ws.Range("J132").Value = 1232.8334
$ws.Range("K132").Value = 7484.8968
$ws.Range("L132").Value = 3698.5002
$ws.Range("M132").Value = -4954.8968
$ws.Range("N132").Value = -8758.5002

$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 755353.1
$ws.Range("I2").Value = 1016.7308
$ws.Range("J2").Value = 2264025.8
$ws.Range("K2").Value = 1016.7308
$ws.Range("L2").Value = 2264025.8
$ws.Range("M2").Value = -903.7308
$ws.Range("N2").Value = -2264251.8
$ws.Range("H37").Value = 8350.625
$ws.Range("I37").Value = 1302.5
$ws.Range("J37").Value = 10700
$ws.Range("K37").Value = 1302.5
$ws.Range("L37").Value = 10700
$ws.Range("M37").Value = -1029.5
$ws.Range("N37").Value = -11246
$ws.Range("H74").Value = 5876.4
$ws.Range("I74").Value = 1113.174
$ws.Range("J74").Value = 21527
$ws.Range("K74").Value = 1113.174
$ws.Range("L74").Value = 21527
$ws.Range("M74").Value = -239.174
$ws.Range("N74").Value = -23275
$ws.Range("H77").Value = 5876.4
$ws.Range("I77").Value = 1113.174
$ws.Range("J77").Value = 21527
$ws.Range("K77").Value = 5565.87
$ws.Range("L77").Value = 107635
$ws.Range("M77").Value = -1197.87
$ws.Range("N77").Value = -116371
$ws.Range("H116").Value = 755353.1
$ws.Range("I116").Value = 1016.7308
$ws.Range("J116").Value = 2264025.8
$ws.Range("K116").Value = 1016.7308
$ws.Range("L116").Value = 2264025.8
$ws.Range("M116").Value = 1277.2692
$ws.Range("N116").Value = -2268613.8
$ws.Range("H124").Value = 9799.799999999999
$ws.Range("J124").Value = 9799.799999999999
$ws.Range("L124").Value = 9799.799999999999
$ws.Range("N124").Value = -19619.8
$ws.Range("H125").Value = 42000
$ws.Range("J125").Value = 42000
$ws.Range("L125").Value = 42000
$ws.Range("N125").Value = -51840

$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 755353.1
$ws.Range("I3").Value = 1016.7308
$ws.Range("J3").Value = 2264025.8
$ws.Range("K3").Value = 1016.7308
$ws.Range("L3").Value = 2264025.8
$ws.Range("M3").Value = -902.7308
$ws.Range("N3").Value = -2264253.8
$ws.Range("H124").Value = 50780
$ws.Range("J124").Value = 50780
$ws.Range("L124").Value = 50780
$ws.Range("N124").Value = -60600

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 542558.8
$ws.Range("I31").Value = 1321.1111
$ws.Range("K31").Value = 1321.1111
$ws.Range("M31").Value = -1026.1111
$ws.Range("H34").Value = 542558.8
$ws.Range("I34").Value = 1321.1111
$ws.Range("K34").Value = 1321.1111
$ws.Range("M34").Value = -1119.1111
$ws.Range("H124").Value = 10520.375
$ws.Range("J124").Value = 10520.375
$ws.Range("L124").Value = 10520.375
$ws.Range("N124").Value = -15430.375
$ws.Range("H132").Value = 1148.017
$ws.Range("I132").Value = 859.54346
$ws.Range("J132").Value = 2168.7693
$ws.Range("K132").Value = 2578.63038
$ws.Range("L132").Value = 6506.3079
$ws.Range("M132").Value = -48.63038000000006
$ws.Range("N132").Value = -11566.3079

$ws = $wb.Worksheets("CUL")
$ws.Range("H113").Value = 525.5192
$ws.Range("I113").Value = 520.7353000000001
$ws.Range("J113").Value = 534.55554
$ws.Range("K113").Value = 1562.2059
$ws.Range("L113").Value = 1603.66662
$ws.Range("M113").Value = 607.7940999999998
$ws.Range("N113").Value = -5943.66662
$ws.Range("H131").Value = 3401.1707
$ws.Range("I131").Value = 5035.385
$ws.Range("J131").Value = 2642.4285
$ws.Range("K131").Value = 15106.155
$ws.Range("L131").Value = 7927.2855
$ws.Range("M131").Value = -10066.155
$ws.Range("N131").Value = -18007.2855

$ws = $wb.Worksheets("GSM")
$ws.Range("H102").Value = 1710.0264
$ws.Range("I102").Value = 1542.9642
$ws.Range("K102").Value = 1542.9642
$ws.Range("M102").Value = 79.03580000000011
$ws.Range("H132").Value = 1718.5873
$ws.Range("I132").Value = 1625.4186
$ws.Range("J132").Value = 1918.9
$ws.Range("K132").Value = 4876.2558
$ws.Range("L132").Value = 5756.700000000001
$ws.Range("M132").Value = -2346.2558
$ws.Range("N132").Value = -10816.7

$ws = $wb.Worksheets("LTW")
$ws.Range("H100").Value = 3937.7
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 5229.5
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 5229.5
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -6311.5

$ws = $wb.Worksheets("WVR")
$ws.Range("H96").Value = 7704038
$ws.Range("I96").Value = 33335334
$ws.Range("J96").Value = 14649.7
$ws.Range("K96").Value = 33335334
$ws.Range("L96").Value = 14649.7
$ws.Range("M96").Value = -33333961
$ws.Range("N96").Value = -17395.7
$ws.Range("H122").Value = 1267.0883
$ws.Range("I122").Value = 1364.625
$ws.Range("J122").Value = 1033
$ws.Range("K122").Value = 4093.875
$ws.Range("L122").Value = 3099
$ws.Range("M122").Value = -1643.875
$ws.Range("N122").Value = -7999
$ws.Range("H132").Value = 3546.638
$ws.Range("I132").Value = 4425.1465
$ws.Range("J132").Value = 1427.8823
$ws.Range("K132").Value = 13275.4395
$ws.Range("L132").Value = 4283.6469
$ws.Range("M132").Value = -10745.4395
$ws.Range("N132").Value = -9343.6469
$ws.Range("H136").Value = 5706.977
$ws.Range("I136").Value = 7567.067
$ws.Range("J136").Value = 1721.0714
$ws.Range("K136").Value = 22701.201
$ws.Range("L136").Value = 5163.2142
$ws.Range("M136").Value = -20151.201
$ws.Range("N136").Value = -10263.2142
